$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2848.1538
$ws.Range("I92").Value = 275.0909
$ws.Range("K92").Value = 275.0909
$ws.Range("M92").Value = 972.9091000000001

$ws.Range("H98").Value = 3844.4722
$ws.Range("I98").Value = 1279.6666
$ws.Range("J98").Value = 16668.5
$ws.Range("K98").Value = 1279.6666
$ws.Range("L98").Value = 16668.5
$ws.Range("M98").Value = 218.3334
$ws.Range("N98").Value = -19664.5

$ws.Range("H122").Value = 3844.4722
$ws.Range("I122").Value = 1279.6666
$ws.Range("J122").Value = 16668.5
$ws.Range("K122").Value = 3838.9998
$ws.Range("L122").Value = 50005.5
$ws.Range("M122").Value = -1388.9998
$ws.Range("N122").Value = -54905.5

$ws.Range("H137").Value = 76462.75
$ws.Range("I137").Value = 60780
$ws.Range("K137").Value = 182340
$ws.Range("M137").Value = -179790

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 1728.3334
$ws.Range("I31").Value = 1728.3334
$ws.Range("K31").Value = 1728.3334
$ws.Range("M31").Value = -1434.3334

$ws.Range("H32").Value = 18437.066
$ws.Range("I32").Value = 18824.78
$ws.Range("K32").Value = 18824.78
$ws.Range("M32").Value = -18537.78

$ws.Range("H61").Value = 8933.467000000001
$ws.Range("I61").Value = 1909.7273
$ws.Range("J61").Value = 28248.75
$ws.Range("K61").Value = 1909.7273
$ws.Range("L61").Value = 28248.75
$ws.Range("M61").Value = -1697.7273
$ws.Range("N61").Value = -28672.75

$ws.Range("H104").Value = 80000
$ws.Range("J104").Value = 80000
$ws.Range("L104").Value = 80000
$ws.Range("N104").Value = -86988

$ws.Range("H132").Value = 1872.8889
$ws.Range("I132").Value = 1513.25
$ws.Range("J132").Value = 4750
$ws.Range("K132").Value = 4539.75
$ws.Range("L132").Value = 14250
$ws.Range("M132").Value = -2009.75
$ws.Range("N132").Value = -19310

$ws.Range("H136").Value = 8933.467000000001
$ws.Range("I136").Value = 1909.7273
$ws.Range("J136").Value = 28248.75
$ws.Range("K136").Value = 5729.1819
$ws.Range("L136").Value = 84746.25
$ws.Range("M136").Value = -3179.1819
$ws.Range("N136").Value = -89846.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9639.041999999999
$ws.Range("I20").Value = 11165.3
$ws.Range("K20").Value = 11165.3
$ws.Range("M20").Value = -10918.3

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1185

$ws.Range("H31").Value = 7695422
$ws.Range("I31").Value = 8336632.5
$ws.Range("J31").Value = 900
$ws.Range("K31").Value = 8336632.5
$ws.Range("L31").Value = 900
$ws.Range("M31").Value = -8336337.5
$ws.Range("N31").Value = -1490

$ws.Range("H34").Value = 7695422
$ws.Range("I34").Value = 8336632.5
$ws.Range("J34").Value = 900
$ws.Range("K34").Value = 8336632.5
$ws.Range("L34").Value = 900
$ws.Range("M34").Value = -8336430.5
$ws.Range("N34").Value = -1304

$ws.Range("H107").Value = 470.75757
$ws.Range("I107").Value = 329
$ws.Range("J107").Value = 848.7778
$ws.Range("K107").Value = 329
$ws.Range("L107").Value = 848.7778
$ws.Range("M107").Value = 1591
$ws.Range("N107").Value = -4688.7778

$ws.Range("H132").Value = 112300.78
$ws.Range("I132").Value = 126225.875
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 378677.625
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -376147.625
$ws.Range("N132").Value = -7760

$ws.Range("H134").Value = 2448.88
$ws.Range("I134").Value = 1686.4706
$ws.Range("K134").Value = 5059.4118
$ws.Range("M134").Value = -2524.4118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 709.6
$ws.Range("J5").Value = 688
$ws.Range("L5").Value = 2064
$ws.Range("N5").Value = -2288

$ws.Range("H131").Value = 429437.1
$ws.Range("I131").Value = 536047.25
$ws.Range("K131").Value = 1608141.75
$ws.Range("M131").Value = -1603101.75

$ws.Range("H132").Value = 1817.4546
$ws.Range("J132").Value = 2082.6667
$ws.Range("L132").Value = 18744.0003
$ws.Range("N132").Value = -23804.0003

$ws.Range("H134").Value = 369.57144
$ws.Range("I134").Value = 369.57144
$ws.Range("K134").Value = 1108.71432
$ws.Range("M134").Value = 3961.28568

$ws.Range("H135").Value = 709.6
$ws.Range("J135").Value = 688
$ws.Range("L135").Value = 6192
$ws.Range("N135").Value = -11262

$ws.Range("H137").Value = 4247.3335
$ws.Range("I137").Value = 4450
$ws.Range("K137").Value = 13350
$ws.Range("M137").Value = -8250

$ws.Range("H139").Value = 5505.0625
$ws.Range("I139").Value = 5505.0625
$ws.Range("K139").Value = 16515.1875
$ws.Range("M139").Value = -11375.1875

$ws.Range("H140").Value = 2549.1667
$ws.Range("I140").Value = 2759
$ws.Range("K140").Value = 8277
$ws.Range("M140").Value = -3097

$ws.Range("H141").Value = 7746
$ws.Range("I141").Value = 6992
$ws.Range("K141").Value = 20976
$ws.Range("M141").Value = -15796

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3114.1875
$ws.Range("I122").Value = 2826.4634
$ws.Range("J122").Value = 4799.4287
$ws.Range("K122").Value = 8479.3902
$ws.Range("L122").Value = 14398.2861
$ws.Range("M122").Value = -6029.3902
$ws.Range("N122").Value = -19298.2861

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3239.842
$ws.Range("I132").Value = 2949.5454
$ws.Range("J132").Value = 3639
$ws.Range("K132").Value = 8848.636200000001
$ws.Range("L132").Value = 10917
$ws.Range("M132").Value = -6318.636200000001
$ws.Range("N132").Value = -15977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 790
$ws.Range("I93").Value = 796.6818
$ws.Range("K93").Value = 796.6818
$ws.Range("M93").Value = 451.3182

$ws.Range("H122").Value = 3029.44
$ws.Range("I122").Value = 2945.1304
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 8835.3912
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -6385.3912
$ws.Range("N122").Value = -16897

$ws.Range("H136").Value = 4539.95
$ws.Range("I136").Value = 3911.5386
$ws.Range("J136").Value = 5707
$ws.Range("K136").Value = 11734.6158
$ws.Range("L136").Value = 17121
$ws.Range("M136").Value = -9184.6158
$ws.Range("N136").Value = -22221

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 881.125
$ws.Range("I107").Value = 884.2222
$ws.Range("K107").Value = 2652.6666
$ws.Range("M107").Value = -732.6666

$ws.Range("H122").Value = 88523.74000000001
$ws.Range("I122").Value = 99388.164
$ws.Range("J122").Value = 1608.3334
$ws.Range("K122").Value = 298164.492
$ws.Range("L122").Value = 4825.0002
$ws.Range("M122").Value = -295714.492
$ws.Range("N122").Value = -9725.0002

$ws.Range("H132").Value = 37245.55
$ws.Range("I132").Value = 56277.23
$ws.Range("J132").Value = 1901
$ws.Range("K132").Value = 168831.69
$ws.Range("L132").Value = 5703
$ws.Range("M132").Value = -168831.69
$ws.Range("N132").Value = -10763

Write-Output "Applied all cell updates."
